$wb = $excel.ActiveWorkbook

# --- AMSIN sheet: add rows 94 and 95 ---
$ws = $wb.Worksheets.Item("AMSIN")

$ws.Range("A94").Value = "2023-02-16"
$ws.Range("B94").Value = 44973.76391283565
$ws.Range("C94").Value = "173cyclefst"
$ws.Range("D94").Value = 269
$ws.Range("E94").Value = 269
$ws.Range("F94").Value = 0
$ws.Range("G94").Value = 4.7
$ws.Range("A94:G94").Style = $ws.Range("A93:G93").Style

$ws.Range("A95").Value = "2023-02-20"
$ws.Range("B95").Value = 44977.39481481481
$ws.Range("C95").Value = "173fnlrun"
$ws.Range("D95").Value = 269
$ws.Range("E95").Value = 269
$ws.Range("F95").Value = 0
$ws.Range("G95").Value = 4.6
$ws.Range("A95:G95").Style = $ws.Range("A93:G93").Style

# --- BETA sheet: add row 34 ---
$ws = $wb.Worksheets.Item("BETA")

$ws.Range("A34").Value = "2023-02-20"
$ws.Range("B34").Value = 44977.57882420139
$ws.Range("C34").Value = "173beta"
$ws.Range("D34").Value = 269
$ws.Range("E34").Value = 267
$ws.Range("F34").Value = 2
$ws.Range("G34").Value = 4.06
$ws.Range("A34:G34").Style = $ws.Range("A33:G33").Style

# --- AMS sheet: fix row 69 style/value, add row 70 ---
$ws = $wb.Worksheets.Item("AMS")

$ws.Range("A69:G69").Style = $ws.Range("A68:G68").Style
$ws.Range("B69").Value = 44963.75031221065

$ws.Range("A70").Value = "2023-02-20"
$ws.Range("B70").Value = 44977.83047720879
$ws.Range("C70").Value = "live173"
$ws.Range("D70").Value = 269
$ws.Range("E70").Value = 268
$ws.Range("F70").Value = 1
$ws.Range("G70").Value = 4.05
